$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 22.65656533333333
$ws.Range("H2").Value = 67.969696
$ws.Range("I2").Value = 0.9268638682343595
$ws.Range("J2").Value = 0.9268638682343595
$ws.Range("M2").Value = 83.91225566666667
$ws.Range("N2").Value = 251.736767
$ws.Range("O2").Value = 0.9556261553553385
$ws.Range("P2").Value = 0.9556261553553385
$ws.Range("Q2").Value = 1901.163502779204
$ws.Range("R2").Value = 17110.47152501283
$ws.Range("S2").Value = 0.885735354938578
$ws.Range("T2").Value = 0.885735354938578
$ws.Range("G3").Value = 22.65656533333333
$ws.Range("H3").Value = 67.969696
$ws.Range("I3").Value = 0.9268638682343595
$ws.Range("J3").Value = 0.9268638682343595
$ws.Range("O3").Value = 0.00439999103960854
$ws.Range("P3").Value = 0.00439999103960854
$ws.Range("Q3").Value = 8.753530164679113
$ws.Range("R3").Value = 78.78177148211201
$ws.Range("S3").Value = 0.004078192715168093
$ws.Range("T3").Value = 0.004078192715168093
$ws.Range("G4").Value = 22.65656533333333
$ws.Range("H4").Value = 67.969696
$ws.Range("I4").Value = 0.9268638682343595
$ws.Range("J4").Value = 0.9268638682343595
$ws.Range("M4").Value = 3.510050666666667
$ws.Range("N4").Value = 10.530152
$ws.Range("O4").Value = 0.03997385360505296
$ws.Range("P4").Value = 0.03997385360505297
$ws.Range("Q4").Value = 79.52569225264357
$ws.Range("R4").Value = 715.7312302737921
$ws.Range("S4").Value = 0.03705032058061339
$ws.Range("T4").Value = 0.0370503205806134
$ws.Range("I5").Value = 0.04016781697437198
$ws.Range("J5").Value = 0.04016781697437198
$ws.Range("M5").Value = 83.91225566666667
$ws.Range("N5").Value = 251.736767
$ws.Range("O5").Value = 0.9556261553553385
$ws.Range("P5").Value = 0.9556261553553385
$ws.Range("Q5").Value = 82.39137400346023
$ws.Range("R5").Value = 741.522366031142
$ws.Range("S5").Value = 0.038385416504236
$ws.Range("T5").Value = 0.038385416504236
$ws.Range("I6").Value = 0.04016781697437198
$ws.Range("J6").Value = 0.04016781697437198
$ws.Range("O6").Value = 0.00439999103960854
$ws.Range("P6").Value = 0.00439999103960854
$ws.Range("S6").Value = 0.0001767380347678725
$ws.Range("T6").Value = 0.0001767380347678725
$ws.Range("I7").Value = 0.04016781697437198
$ws.Range("J7").Value = 0.04016781697437198
$ws.Range("M7").Value = 3.510050666666667
$ws.Range("N7").Value = 10.530152
$ws.Range("O7").Value = 0.03997385360505296
$ws.Range("P7").Value = 0.03997385360505297
$ws.Range("Q7").Value = 3.446432168350222
$ws.Range("R7").Value = 31.017889515152
$ws.Range("S7").Value = 0.001605662435368107
$ws.Range("T7").Value = 0.001605662435368107
$ws.Range("G8").Value = 0.8058883333333333
$ws.Range("H8").Value = 2.417665
$ws.Range("I8").Value = 0.03296831479126849
$ws.Range("J8").Value = 0.03296831479126849
$ws.Range("M8").Value = 83.91225566666667
$ws.Range("N8").Value = 251.736767
$ws.Range("O8").Value = 0.9556261553553385
$ws.Range("P8").Value = 0.9556261553553385
$ws.Range("Q8").Value = 67.62390786545055
$ws.Range("R8").Value = 608.615170789055
$ws.Range("S8").Value = 0.03150538391252445
$ws.Range("T8").Value = 0.03150538391252444
$ws.Range("G9").Value = 0.8058883333333333
$ws.Range("H9").Value = 2.417665
$ws.Range("I9").Value = 0.03296831479126849
$ws.Range("J9").Value = 0.03296831479126849
$ws.Range("O9").Value = 0.00439999103960854
$ws.Range("P9").Value = 0.00439999103960854
$ws.Range("Q9").Value = 0.3113608674311111
$ws.Range("R9").Value = 2.80224780688
$ws.Range("S9").Value = 0.0001450602896725751
$ws.Range("T9").Value = 0.000145060289672575
$ws.Range("G10").Value = 0.8058883333333333
$ws.Range("H10").Value = 2.417665
$ws.Range("I10").Value = 0.03296831479126849
$ws.Range("J10").Value = 0.03296831479126849
$ws.Range("M10").Value = 3.510050666666667
$ws.Range("N10").Value = 10.530152
$ws.Range("O10").Value = 0.03997385360505296
$ws.Range("P10").Value = 0.03997385360505297
$ws.Range("Q10").Value = 2.828708881675555
$ws.Range("R10").Value = 25.45837993508
$ws.Range("S10").Value = 0.001317870589071469
$ws.Range("T10").Value = 0.001317870589071469
